# PRI-22079 Export Campaign to Excel - All Tabs - Update text patterns in Content Restrictions
# Rename the "Flow Chart template tables" sheet to "Flow Chart Template Tables"
# (title-case the word "Template Tables").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flow Chart template tables")
$ws.Name = "Flow Chart Template Tables"
